$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Shared-string text clean-up: "my reference" -> "my references"
#    and "my drive" -> "my OS-drive". Update every sheet that carries
#    these header cells (fieldnames, URL, comments) so the old shared
#    strings become orphaned and the table stays compact.
# ------------------------------------------------------------------
$headerSheets = @("fieldnames", "URL", "comments")
foreach ($name in $headerSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Value = "my references"
    $ws.Range("C1").Value = "my OS-drive"
}

# ------------------------------------------------------------------
# 2. The "icons" sheet gains the same header row the other sheets
#    already have: A1/B1/C1 = All-time favorites / my references /
#    my OS-drive, with B1:C1 using the "text" style (left/top align)
#    that matches the other header rows.
# ------------------------------------------------------------------
$wsIcons = $wb.Worksheets.Item("icons")
$wsIcons.Range("A1").Value = "All-time favorites"
$wsIcons.Range("B1").Value = "my references"
$wsIcons.Range("C1").Value = "my OS-drive"
$wsIcons.Range("B1:C1").NumberFormat = "@"
$wsIcons.Range("B1:C1").HorizontalAlignment = -4131
$wsIcons.Range("B1:C1").VerticalAlignment = -4160

# ------------------------------------------------------------------
# 3. Selection clean-up on every sheet. "icons" must stay the active
#    sheet/tab, so it is touched last.
# ------------------------------------------------------------------
$wb.Worksheets.Item("fieldnames").Range("A1:C1").Select()
$wb.Worksheets.Item("URL").Range("A1:C1").Select()
$wb.Worksheets.Item("color").Range("A3").Select()
$wb.Worksheets.Item("comments").Range("A1:C1").Select()
$wb.Worksheets.Item("icons").Range("A1:C1").Select()
